{"js": "// \"download articles with pandoc title blocks\"\n//\n// Turns the old hand-rolled title block:\n//   [bookmark] Heading1 paragraph \"On Pilgrimage - October/November 1972\" [/bookmark]\n//   paragraph, bold run \"By Dorothy Day\"\n// into a pandoc-generated-looking title block:\n//   Title-style paragraph, text split word-by-word into separate runs\n//   Authors-style paragraph, text split word-by-word into separate runs (\"By \" dropped)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst authorPara = paragraphs.items[1];\n\n// Best-effort: drop the legacy bookmark that wrapped the old heading paragraph\n// (the \"on-pilgrimage---octobernovember-1972\" slug). Some hosts don't surface\n// document bookmarks for mutation, so tolerate that silently either way.\ntry {\n  context.document.deleteBookmark(\"on-pilgrimage---octobernovember-1972\");\n  await context.sync();\n} catch (e) {\n  // no-op: nothing we can do if the host doesn't expose bookmark deletion\n}\n\nfunction wordRuns(words) {\n  return words\n    .map((w) => `<w:r><w:t xml:space=\"preserve\">${w}</w:t></w:r>`)\n    .join(\"\");\n}\n\nfunction packageXml(paragraphXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    `<w:body>${paragraphXml}</w:body>` +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// --- Title paragraph: \"On Pilgrimage - October/November 1972\", one run per word/space ---\nconst titleWords = [\"On\", \" \", \"Pilgrimage\", \" \", \"-\", \" \", \"October\", \"/\", \"November\", \" \", \"1972\"];\nconst titleParaXml = `<w:p><w:pPr><w:pStyle w:val=\"Title\"/></w:pPr>${wordRuns(titleWords)}</w:p>`;\ntitlePara.getRange(\"Whole\").insertOoxml(packageXml(titleParaXml), \"Replace\");\nawait context.sync();\n\n// Re-fetch paragraphs since the OOXML insert reshapes the paragraph collection.\nparagraphs.load(\"items\");\nawait context.sync();\nconst authorParaFresh = paragraphs.items[1];\n\n// --- Authors paragraph: \"Dorothy Day\" (no \"By \"), one run per word/space ---\nconst authorWords = [\"Dorothy\", \" \", \"Day\"];\nconst authorParaXml = `<w:p><w:pPr><w:pStyle w:val=\"Authors\"/></w:pPr>${wordRuns(authorWords)}</w:p>`;\nauthorParaFresh.getRange(\"Whole\").insertOoxml(packageXml(authorParaXml), \"Replace\");\nawait context.sync();\n", "ps1": "# \"download articles with pandoc title blocks\"\n#\n# Turns the old hand-rolled title block:\n#   [bookmark] Heading1 paragraph \"On Pilgrimage - October/November 1972\" [/bookmark]\n#   paragraph, bold run \"By Dorothy Day\"\n# into a pandoc-generated-looking title block:\n#   Title-style paragraph, text split word-by-word into separate runs\n#   Authors-style paragraph, text split word-by-word into separate runs (\"By \" dropped)\n\n$d = $word.ActiveDocument\n\n# Best-effort: drop the legacy bookmark that wrapped the old heading paragraph\n# (the \"on-pilgrimage---octobernovember-1972\" slug). Some hosts don't surface\n# document bookmarks for mutation, so tolerate that silently either way.\ntry {\n    if ($d.Bookmarks.Exists(\"on-pilgrimage---octobernovember-1972\")) {\n        $d.Bookmarks(\"on-pilgrimage---octobernovember-1972\").Delete()\n    }\n} catch {\n    # no-op: nothing we can do if the host doesn't expose bookmark deletion\n}\n\nfunction Get-WordRunsXml($words) {\n    $sb = \"\"\n    foreach ($w in $words) {\n        $sb += \"<w:r><w:t xml:space=`\"preserve`\">$w</w:t></w:r>\"\n    }\n    return $sb\n}\n\nfunction Get-PackageXml($paragraphXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n        '<pkg:xmlData>' + `\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n        \"<w:body>$paragraphXml</w:body>\" + `\n        '</w:document>' + `\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# --- Title paragraph: \"On Pilgrimage - October/November 1972\", one run per word/space ---\n$titleWords = @(\"On\", \" \", \"Pilgrimage\", \" \", \"-\", \" \", \"October\", \"/\", \"November\", \" \", \"1972\")\n$titleParaXml = \"<w:p><w:pPr><w:pStyle w:val=`\"Title`\"/></w:pPr>\" + (Get-WordRunsXml $titleWords) + \"</w:p>\"\n$p1 = $d.Paragraphs(1)\n$p1.Range.InsertXML((Get-PackageXml $titleParaXml)) | Out-Null\n\n# --- Authors paragraph: \"Dorothy Day\" (no \"By \"), one run per word/space ---\n$authorWords = @(\"Dorothy\", \" \", \"Day\")\n$authorParaXml = \"<w:p><w:pPr><w:pStyle w:val=`\"Authors`\"/></w:pPr>\" + (Get-WordRunsXml $authorWords) + \"</w:p>\"\n$p2 = $d.Paragraphs(2)\n$p2.Range.InsertXML((Get-PackageXml $authorParaXml)) | Out-Null\n"}
